$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the unmodified "data" style (s=3) used to restore
# formatting on cells where assigning a percent-like string would otherwise
# make Excel auto-convert the text into a numeric percentage.
$fmtSource = $ws.Range("C2")

$ws.Range("E2").Value = "2026-02-07 22:48:36"
$ws.Range("E3").Value = "2026-02-07 22:48:38"
$ws.Range("I3").Value = "0.1 mm"
$ws.Range("E4").Value = "2026-02-07 22:48:41"
$ws.Range("J4").Value = "1004.1 hPa"
$ws.Range("N4").Value = "8.6 °C 22:29 TU"
$ws.Range("O4").Value = "11.8 °C"
$ws.Range("E5").Value = "2026-02-07 22:48:44"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "71%"
$fmtSource.Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("O5").Value = "-3.6 °C"
$ws.Range("E6").Value = "2026-02-07 22:48:46"
$ws.Range("E7").Value = "2026-02-07 22:48:49"
$ws.Range("E8").Value = "2026-02-07 22:48:51"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "63%"
$fmtSource.Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("E9").Value = "2026-02-07 22:48:54"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "75%"
$fmtSource.Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("N9").Value = "3.8 °C 22:09 TU"
$ws.Range("O9").Value = "10.6 °C"
$ws.Range("E10").Value = "2026-02-07 22:48:57"
$ws.Range("E11").Value = "2026-02-07 22:48:59"
$ws.Range("E12").Value = "2026-02-07 22:49:02"
$ws.Range("O12").Value = "10.0 °C"
$ws.Range("E13").Value = "2026-02-07 22:49:05"
$ws.Range("E14").Value = "2026-02-07 22:49:08"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "61%"
$fmtSource.Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("O14").Value = "11.7 °C"
$ws.Range("E15").Value = "2026-02-07 22:49:10"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "71%"
$fmtSource.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("N15").Value = "4.3 °C 22:05 TU"
$ws.Range("O15").Value = "10.2 °C"
$ws.Range("E16").Value = "2026-02-07 22:49:13"
$ws.Range("I16").Value = "0.7 mm"
$ws.Range("O16").Value = "-5.8 °C"
$ws.Range("E17").Value = "2026-02-07 22:49:16"
$ws.Range("L17").Value = "74.9 km/h - 230º 22:01 TU"
$ws.Range("E18").Value = "2026-02-07 22:49:18"
$ws.Range("E19").Value = "2026-02-07 22:49:21"
$ws.Range("E20").Value = "2026-02-07 22:49:24"
$ws.Range("I20").Value = "2.7 mm"
$ws.Range("L20").Value = "42.8 km/h - 193º 22:24 TU"
$ws.Range("E21").Value = "2026-02-07 22:49:26"
$ws.Range("E22").Value = "2026-02-07 22:49:29"
$ws.Range("I22").Value = "2.4 mm"
$ws.Range("E23").Value = "2026-02-07 22:49:32"
$ws.Range("E24").Value = "2026-02-07 22:49:34"
$ws.Range("I24").Value = "1.3 mm"
$ws.Range("J24").Value = "1007.1 hPa"
$ws.Range("E25").Value = "2026-02-07 22:49:37"
$ws.Range("L25").Value = "24.8 km/h - 227º 22:10 TU"
$ws.Range("E26").Value = "2026-02-07 22:49:40"
$ws.Range("L26").Value = "42.8 km/h - 218º 22:19 TU"
$ws.Range("E27").Value = "2026-02-07 22:49:43"
$ws.Range("I27").Value = "2.5 mm"
$ws.Range("E28").Value = "2026-02-07 22:49:45"
$ws.Range("J28").Value = "1004.2 hPa"
$ws.Range("L28").Value = "42.1 km/h - 254º 22:21 TU"
$ws.Range("E29").Value = "2026-02-07 22:49:48"
$ws.Range("E30").Value = "2026-02-07 22:49:50"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "81%"
$fmtSource.Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("E31").Value = "2026-02-07 22:49:53"
$ws.Range("E32").Value = "2026-02-07 22:49:55"
$ws.Range("I32").Value = "1.0 mm"
$ws.Range("O32").Value = "4.4 °C"
$ws.Range("E33").Value = "2026-02-07 22:49:58"
$ws.Range("E34").Value = "2026-02-07 22:50:01"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "68%"
$fmtSource.Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("L34").Value = "41.0 km/h - 194º 22:09 TU"
$ws.Range("E35").Value = "2026-02-07 22:50:03"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "69%"
$fmtSource.Copy()
$ws.Range("H35").PasteSpecial(-4122)
$ws.Range("J35").Value = "1007.3 hPa"
$ws.Range("E36").Value = "2026-02-07 22:50:05"
$ws.Range("E37").Value = "2026-02-07 22:50:08"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "80%"
$fmtSource.Copy()
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("J37").Value = "1006.0 hPa"
$ws.Range("O37").Value = "4.7 °C"
$ws.Range("E38").Value = "2026-02-07 22:50:11"
$ws.Range("N38").Value = "7.4 °C 22:26 TU"
$ws.Range("O38").Value = "11.9 °C"
$ws.Range("E39").Value = "2026-02-07 22:50:13"
$ws.Range("E40").Value = "2026-02-07 22:50:16"
$ws.Range("O40").Value = "5.3 °C"
$ws.Range("E41").Value = "2026-02-07 22:50:19"
$ws.Range("J41").Value = "1006.3 hPa"
$ws.Range("E42").Value = "2026-02-07 22:50:21"
$ws.Range("E43").Value = "2026-02-07 22:50:24"
$ws.Range("E44").Value = "2026-02-07 22:50:27"
$ws.Range("E45").Value = "2026-02-07 22:50:29"
$ws.Range("E46").Value = "2026-02-07 22:50:32"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "74%"
$fmtSource.Copy()
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("J46").Value = "1007.4 hPa"
